# Updates the "想去人数" (want-to-go count) values in column F across the
# "展览" and "全部类型" sheets (and one row in "演出"), reflecting a fresh
# data pull (gh-pages output regenerated at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1791   # was 1790
$ws1.Range("F4").Value  = 458    # was 454
$ws1.Range("F5").Value  = 81     # was 80
$ws1.Range("F8").Value  = 339    # was 337
$ws1.Range("F9").Value  = 1740   # was 1738
$ws1.Range("F11").Value = 1424   # was 1423
$ws1.Range("F12").Value = 814    # was 811
$ws1.Range("F14").Value = 683    # was 682
$ws1.Range("F15").Value = 12829  # was 12822
$ws1.Range("F16").Value = 12820  # was 12813
$ws1.Range("F23").Value = 2011   # was 2010
$ws1.Range("F24").Value = 32     # was 31
$ws1.Range("F25").Value = 6      # was 5
$ws1.Range("F27").Value = 50     # was 43
$ws1.Range("F28").Value = 252    # was 251

# --- Sheet "演出" (sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F10").Value = 79     # was 78

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 1791   # was 1790
$ws4.Range("F6").Value  = 458    # was 454
$ws4.Range("F8").Value  = 81     # was 80
$ws4.Range("F13").Value = 339    # was 337
$ws4.Range("F14").Value = 1740   # was 1738
$ws4.Range("F16").Value = 1424   # was 1423
$ws4.Range("F17").Value = 814    # was 811
$ws4.Range("F20").Value = 683    # was 682
$ws4.Range("F21").Value = 12829  # was 12822
$ws4.Range("F22").Value = 12820  # was 12813
$ws4.Range("F31").Value = 2011   # was 2010
$ws4.Range("F32").Value = 32     # was 31
$ws4.Range("F33").Value = 6      # was 5
$ws4.Range("F37").Value = 50     # was 43
$ws4.Range("F38").Value = 252    # was 251
$ws4.Range("F40").Value = 79     # was 78
